$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Euro" row symbol: EUR -> EUR€
$ws.Range("B3").Value = "EUR€"

# New cell H1 with the bare euro sign, using an explicit (non-themed) Calibri font
$ws.Range("H1").Value = "€"
$ws.Range("H1").Font.Name = "Calibri"

# Update the current selection to C3
$ws.Range("C3").Select() | Out-Null
